# Fix Training Data Issue (#48)
# Data was taken from 1 day off due to way NBA stats were shown.
# The "Date" column (BF) on Sheet1 holds the literal string "5-13-2007-08"
# for every team row (rows 2-31); it should read "2008-05-13" instead.
#
# NOTE: "2008-05-13" looks like an ISO date to Excel's smart-entry parser,
# so a naive .Value assignment would silently convert the cell to a real
# date serial number (and pull in a new number-format/style). We force the
# cell to Text format first so the value is kept as a literal string (same
# as the original "inlineStr"/text cell), then restore the cell style so
# no extra formatting is left behind on the cells themselves.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rng = $ws.Range("BF2:BF31")
$rng.NumberFormat = "@"
$rng.Value = "2008-05-13"
$rng.Style = "Normal"
